# Kilimanjaro Weekly Scoreboard - append this week's workout rows (205-207)
# for 2024-07-13 (Excel serial date 45486), matching the rows already
# logged for previous days/participants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the date style (numFmtId 14, short date) from the last existing
# date cell onto the three new date cells before writing their values, so
# the new <c> elements get the same s="1" style index as the rest of
# column B instead of the engine minting a brand-new style.
$ws.Range("B204").Copy($ws.Range("B205:B207"))

$rows = @(
    @{ Row = 205; Participant = "Jeremiah"; Date = 45486; Workout = "Workout"; Duration = 50;  Distance = 0;    Elevation = 0;   Z1 = 37; Z2 = 13; Z3 = 0; Z4 = 0; Z5 = 0; Level = "Wily Hyena";    Week = 5 },
    @{ Row = 206; Participant = "Jeremiah"; Date = 45486; Workout = "Run";     Duration = 23;  Distance = 2.71; Elevation = 173; Z1 = 0;  Z2 = 5;  Z3 = 9; Z4 = 5; Z5 = 0; Level = "Wily Hyena";    Week = 5 },
    @{ Row = 207; Participant = "Steven";   Date = 45486; Workout = "Walk";    Duration = 29;  Distance = 1.61; Elevation = 26;  Z1 = 29; Z2 = 0;  Z3 = 0; Z4 = 0; Z5 = 0; Level = "Brave Leopard"; Week = 5 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Participant
    $ws.Cells.Item($row, 2).Value = $r.Date
    $ws.Cells.Item($row, 3).Value = $r.Workout
    $ws.Cells.Item($row, 4).Value = $r.Duration
    $ws.Cells.Item($row, 5).Value = $r.Distance
    $ws.Cells.Item($row, 6).Value = $r.Elevation
    $ws.Cells.Item($row, 7).Value = $r.Z1
    $ws.Cells.Item($row, 8).Value = $r.Z2
    $ws.Cells.Item($row, 9).Value = $r.Z3
    $ws.Cells.Item($row, 10).Value = $r.Z4
    $ws.Cells.Item($row, 11).Value = $r.Z5
    $ws.Cells.Item($row, 12).Value = $r.Level
    $ws.Cells.Item($row, 13).Value = $r.Week
}

# Reproduce the author's final view state: the frozen pane scrolled down a
# few rows and the cursor resting just below the newly-added data.
$ws.Range("A208").Select()
